$p = $ppt.ActivePresentation
$s = $p.Slides.Item(20)

# Mark the slide as hidden during slide show (adds show="0" to <p:sld>)
$s.SlideShowTransition.Hidden = $true

# Give the slide a slow transition lasting 2 seconds
# (Speed -> spd="slow" on <p:transition>, Duration -> p14:dur="2000")
$s.SlideShowTransition.Speed = 1  # ppTransitionSpeedSlow
$s.SlideShowTransition.Duration = 2
